$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new job posting row (row 6) with Job_Id = JD_005
$ws.Range("A6").Value = "JD_005"
$ws.Range("B6").Value = "Senior Devops Engineer"
$ws.Range("C6").Value = "Demo"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 4
